# Updated some advice on readme
# Add four new rows of model notes to the bottom of the table and
# tidy up the sheet view (turn off "show formulas", move the
# selection to the new last cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content -----------------------------------------------------
$ws.Range("B4").Value = "Mobile Net pretrained model"
$ws.Range("C4").Value = "Better than most models "
$ws.Range("B5").Value = "change input shape "
$ws.Range("B6").Value = "try removing rescaling "
$ws.Range("B7").Value = "greyscale "

# Match the row height used by the rest of the wrapped-text rows.
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17

# --- View tweaks -------------------------------------------------------
# Turn off "Show Formulas" for the sheet.
$excel.ActiveWindow.DisplayFormulas = $false

# Move the active selection to the last cell that was filled in.
$ws.Range("B7").Select() | Out-Null
